$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric while preserving the original
# "text" cell type (the source sheet stores every data cell, even price/
# hour columns, as inline text rather than numbers). Temporarily switching
# the cell to a Text number format forces Excel to keep the literal string
# instead of auto-converting it to a Number, then we restore the cell style
# so no stray formatting/style record is left behind.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "243.47"
Set-TextValue $ws.Range("G2") "15"
Set-TextValue $ws.Range("D3") "24.05"
Set-TextValue $ws.Range("G3") "15"
Set-TextValue $ws.Range("D4") "5.402"
Set-TextValue $ws.Range("G4") "15"
Set-TextValue $ws.Range("D5") "0.05890"
Set-TextValue $ws.Range("G5") "15"
Set-TextValue $ws.Range("D6") "3.393"
Set-TextValue $ws.Range("G6") "15"
Set-TextValue $ws.Range("D7") "6.506"
Set-TextValue $ws.Range("G7") "15"
Set-TextValue $ws.Range("D8") "0.8113"
Set-TextValue $ws.Range("G8") "15"
Set-TextValue $ws.Range("D9") "0.9235"
Set-TextValue $ws.Range("G9") "15"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D10") "0.01094"
$ws.Range("E10").Value = "9OneONEBestin24h"
Set-TextValue $ws.Range("G10") "15"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1419"
$ws.Range("E11").Value = "10WazirXWRX"
Set-TextValue $ws.Range("G11") "15"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07404"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
Set-TextValue $ws.Range("G12") "15"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D13") "0.03072"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws.Range("G13") "15"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03056"
$ws.Range("E14").Value = "13BitrueCoinBTR"
Set-TextValue $ws.Range("G14") "15"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09355"
$ws.Range("E15").Value = "14BitMartTokenBMX"
Set-TextValue $ws.Range("G15") "15"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D16") "3.851"
$ws.Range("E16").Value = "15MCDexMCB"
Set-TextValue $ws.Range("G16") "15"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001574"
$ws.Range("E17").Value = "16BitForexTokenBF"
Set-TextValue $ws.Range("G17") "15"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04709"
$ws.Range("E18").Value = "17CoinExTokenCET"
Set-TextValue $ws.Range("G18") "15"
Set-TextValue $ws.Range("D19") "0.005941"
Set-TextValue $ws.Range("G19") "15"
Set-TextValue $ws.Range("D20") "0.001249"
Set-TextValue $ws.Range("G20") "15"
Set-TextValue $ws.Range("D21") "0.004727"
Set-TextValue $ws.Range("G21") "15"
Set-TextValue $ws.Range("D22") "0.00008808"
Set-TextValue $ws.Range("G22") "15"
Set-TextValue $ws.Range("D23") "3.554"
Set-TextValue $ws.Range("G23") "15"
Set-TextValue $ws.Range("G24") "15"
Set-TextValue $ws.Range("D25") "0.3224"
Set-TextValue $ws.Range("G25") "15"
Set-TextValue $ws.Range("D26") "0.1331"
Set-TextValue $ws.Range("G26") "15"
Set-TextValue $ws.Range("D27") "0.0002655"
$ws.Range("E27").Value = "26UpBotsUBXT"
Set-TextValue $ws.Range("G27") "15"
Set-TextValue $ws.Range("G28") "15"
Set-TextValue $ws.Range("G29") "15"
Set-TextValue $ws.Range("G30") "15"
Set-TextValue $ws.Range("G31") "15"
Set-TextValue $ws.Range("G32") "15"
Set-TextValue $ws.Range("G33") "15"
Set-TextValue $ws.Range("G34") "15"
Set-TextValue $ws.Range("G35") "15"
Set-TextValue $ws.Range("G36") "15"
Set-TextValue $ws.Range("G37") "15"
Set-TextValue $ws.Range("G38") "15"
Set-TextValue $ws.Range("G39") "15"
Set-TextValue $ws.Range("D40") "0.03880"
Set-TextValue $ws.Range("G40") "15"
Set-TextValue $ws.Range("D41") "0.006345"
Set-TextValue $ws.Range("G41") "15"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.003503"
$ws.Range("E42").Value = "41CEJICEJI"
Set-TextValue $ws.Range("G42") "15"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D43") "0.1065"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue $ws.Range("G43") "15"
Set-TextValue $ws.Range("D44") "0.008553"
Set-TextValue $ws.Range("G44") "15"
Set-TextValue $ws.Range("D45") "0.00005219"
Set-TextValue $ws.Range("G45") "15"
Set-TextValue $ws.Range("G46") "15"
Set-TextValue $ws.Range("D47") "0.6716"
Set-TextValue $ws.Range("G47") "15"
Set-TextValue $ws.Range("D48") "0.001943"
Set-TextValue $ws.Range("G48") "15"
Set-TextValue $ws.Range("G49") "15"
Set-TextValue $ws.Range("G50") "15"
Set-TextValue $ws.Range("G51") "15"
